$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous row 3 (A3="Q_to_net1", B3=0) moves down to row 4, and a new row 3
# ("param_Q_to_demand2" / "Q_net1_demand2") is inserted above it.

# Give the new row 4's label cell the same style as the existing label cells (A2/A3),
# then write its content.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A4").Value = "Q_to_net1"
$ws.Range("B4").Value = 0

# Overwrite the old row 3 with the newly inserted content.
$ws.Range("A3").Value = "param_Q_to_demand2"
$ws.Range("B3").Value = "Q_net1_demand2"
